$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-03-24 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-03-25 Monday", 2) | Out-Null
$d.Content.Find.Execute("95×57=5415", $true, $false, $false, $false, $false, $true, 1, $false, "69×78=5382", 2) | Out-Null
$d.Content.Find.Execute("61×23=1403", $true, $false, $false, $false, $false, $true, 1, $false, "65×54=3510", 2) | Out-Null
$d.Content.Find.Execute("35×81=2835", $true, $false, $false, $false, $false, $true, 1, $false, "57×58=3306", 2) | Out-Null
$d.Content.Find.Execute("21×25=525", $true, $false, $false, $false, $false, $true, 1, $false, "66×77=5082", 2) | Out-Null
$d.Content.Find.Execute("81×31=2511", $true, $false, $false, $false, $false, $true, 1, $false, "50×72=3600", 2) | Out-Null
$d.Content.Find.Execute("62×40=2480", $true, $false, $false, $false, $false, $true, 1, $false, "11×83=913", 2) | Out-Null
$d.Content.Find.Execute("31×12=372", $true, $false, $false, $false, $false, $true, 1, $false, "64×79=5056", 2) | Out-Null
$d.Content.Find.Execute("34×34=1156", $true, $false, $false, $false, $false, $true, 1, $false, "53×14=742", 2) | Out-Null
$d.Content.Find.Execute("20×63=1260", $true, $false, $false, $false, $false, $true, 1, $false, "42×65=2730", 2) | Out-Null
$d.Content.Find.Execute("60×17=1020", $true, $false, $false, $false, $false, $true, 1, $false, "69×48=3312", 2) | Out-Null
$d.Content.Find.Execute("76×58=4408", $true, $false, $false, $false, $false, $true, 1, $false, "57×67=3819", 2) | Out-Null
$d.Content.Find.Execute("66×67=4422", $true, $false, $false, $false, $false, $true, 1, $false, "98×59=5782", 2) | Out-Null
$d.Content.Find.Execute("62×84=5208", $true, $false, $false, $false, $false, $true, 1, $false, "43×43=1849", 2) | Out-Null
$d.Content.Find.Execute("45×74=3330", $true, $false, $false, $false, $false, $true, 1, $false, "59×61=3599", 2) | Out-Null
$d.Content.Find.Execute("36×95=3420", $true, $false, $false, $false, $false, $true, 1, $false, "99×67=6633", 2) | Out-Null
$d.Content.Find.Execute("49×96=4704", $true, $false, $false, $false, $false, $true, 1, $false, "54×66=3564", 2) | Out-Null
$d.Content.Find.Execute("44×52=2288", $true, $false, $false, $false, $false, $true, 1, $false, "39×96=3744", 2) | Out-Null
$d.Content.Find.Execute("56×17=952", $true, $false, $false, $false, $false, $true, 1, $false, "41×63=2583", 2) | Out-Null
$d.Content.Find.Execute("52×76=3952", $true, $false, $false, $false, $false, $true, 1, $false, "52×62=3224", 2) | Out-Null
$d.Content.Find.Execute("49×60=2940", $true, $false, $false, $false, $false, $true, 1, $false, "68×65=4420", 2) | Out-Null
$d.Content.Find.Execute("89×42=3738", $true, $false, $false, $false, $false, $true, 1, $false, "68×43=2924", 2) | Out-Null
$d.Content.Find.Execute("89×11=979", $true, $false, $false, $false, $false, $true, 1, $false, "13×91=1183", 2) | Out-Null
$d.Content.Find.Execute("84×19=1596", $true, $false, $false, $false, $false, $true, 1, $false, "82×82=6724", 2) | Out-Null
$d.Content.Find.Execute("96×90=8640", $true, $false, $false, $false, $false, $true, 1, $false, "33×43=1419", 2) | Out-Null
$d.Content.Find.Execute("39×18=702", $true, $false, $false, $false, $false, $true, 1, $false, "62×93=5766", 2) | Out-Null
